# DVM - Business Rules.xlsx
# Add a new column I "Test Case Exists? (if applicable)" right after the
# "Example (optional)" column (H), indicating whether a DVM test case
# exists for the given rule / error code. All the former columns I:O
# (Related Rules .. Retired) shift right to J:P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new column, shifting I:O -> J:P -------------------------
$ws.Columns("I:I").Insert()

# --- 2. Header cell for the new column --------------------------------------
$hdr = $ws.Cells.Item(1, 9)
$hdr.Value = "Test Case Exists? (if applicable)"
# Bold header to match the rest of row 1, no wrap/fill (matches the other
# header cells' look-and-feel).
$hdr.Font.Bold = $true
$hdr.WrapText = $false

# --- 3. Fill in the new column's values -------------------------------------
# Rows 2-13   -> Scope "DVM DB"              -> not applicable
# Rows 14-27  -> Scope "DVM Configuration QC" -> test cases exist
# Rows 28-65  -> Scope "DVM Processing Errors" -> mixed yes/no per rule

$naRows   = @(2, 13)
$yesRanges = @(
    @(14, 27),
    @(29, 32),
    @(34, 34),
    @(38, 39),
    @(48, 49),
    @(52, 53),
    @(62, 64)
)
$noRanges = @(
    @(28, 28),
    @(33, 33),
    @(35, 37),
    @(40, 47),
    @(50, 51),
    @(54, 61),
    @(65, 65)
)

for ($r = $naRows[0]; $r -le $naRows[1]; $r++) {
    $ws.Cells.Item($r, 9).Value = "N/A"
}

# "no" is written before "yes" so the shared-strings table lists them in
# the same order as the source workbook (no, then yes).
foreach ($range in $noRanges) {
    for ($r = $range[0]; $r -le $range[1]; $r++) {
        $ws.Cells.Item($r, 9).Value = "no"
    }
}

foreach ($range in $yesRanges) {
    for ($r = $range[0]; $r -le $range[1]; $r++) {
        $ws.Cells.Item($r, 9).Value = "yes"
    }
}

# --- 4. Restore the view: scrolled back to top, new column selected --------
$ws.Range("I4").Select()
